# Release-Notes.xlsx update:
#  - A new folder entry "Azure Virtual Machine And Compute" was added.
#  - The folder "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals" was
#    refreshed (new "Last Updated" timestamp) and re-sorted to the top.
#  - All other "Folder Inventory" rows shift down to make room, and the
#    stale "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals" row (with
#    its old timestamp) is dropped since it is superseded by the new one.
#  - Metadata / Summary sheets are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folder Inventory")

# --- Step 1: shift rows 9..71 down to 10..72 (bottom-up so nothing is clobbered) ---
for ($r = 71; $r -ge 9; $r--) {
    $dstRow = $r + 1
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($dstRow, $c).Value = $ws.Cells.Item($r, $c).Value()
    }
}

# --- Step 2: shift rows 2..7 down to 4..9 (bottom-up) ---
for ($r = 7; $r -ge 2; $r--) {
    $dstRow = $r + 2
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($dstRow, $c).Value = $ws.Cells.Item($r, $c).Value()
    }
}

# --- Step 3: write the two new top rows ---
$ws.Cells.Item(2, 1).Value = "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals"
$ws.Cells.Item(2, 2).Value = "Azure_AI_Foundry_and_Semantic_Kernel_Fundamentals"
$ws.Cells.Item(2, 3).Value = "2025-06-12 15:19:27 +0530"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "Root"

$ws.Cells.Item(3, 1).Value = "Azure Virtual Machine And Compute"
$ws.Cells.Item(3, 2).Value = "Azure Virtual Machine And Compute"
$ws.Cells.Item(3, 3).Value = "2025-06-12 15:19:27 +0530"
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = "Root"

# --- Step 4: refresh the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value = "2025-06-12 09:49:46 UTC"
$meta.Cells.Item(4, 2).Value = 71
$meta.Cells.Item(5, 2).Value = "9"

# --- Step 5: refresh the "Summary" sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(2, 2).Value = 71
$summary.Cells.Item(3, 2).Value = 71
$summary.Cells.Item(5, 2).Value = "2025-06-12 15:19:27 +0530"
